{"js": "// Fixed alert confirmation popup not working in DELETE\n//\n// 1) Row 2's first paragraph currently has its Greek text split across two\n//    runs (with a leftover <w:bookmarkStart/bookmarkEnd name=\"_GoBack\">\n//    sitting between them, left over from where the cursor last was).\n//    Merge that back into a single contiguous run and drop the stray\n//    bookmark.\n// 2) Row 3's content cell is empty. Fill it in with the new feedback\n//    (that there was no confirmation on delete, so it could be clicked by\n//    accident) and the corresponding fix note (an alert/confirm message\n//    was added on delete) - moving the \"_GoBack\" bookmark to sit at the very\n//    end of that new content, which is where the author's cursor ended up.\n\nconst OOXML_NS =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" ' +\n  'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>{BODY}</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapOoxml(bodyXml) {\n  return OOXML_NS.replace(\"{BODY}\", bodyXml);\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// --- Part 1: merge the split run in row 2 (index 1) and drop the bookmark\nconst row2Cell = table.getCell(1, 1);\nconst row2Paragraphs = row2Cell.body.paragraphs;\nrow2Paragraphs.load(\"text\");\nawait context.sync();\n\nconst mergedText =\n  \"\u0395\u03bd\u03ce \u03ad\u03c7\u03b5\u03b9 \u03bd\u03cc\u03b7\u03bc\u03b1 \u03bd\u03b1 \u03ad\u03c7\u03bf\u03c5\u03bc\u03b5 \u03be\u03b5\u03c7\u03c9\u03c1\u03b9\u03c3\u03c4\u03ae \u03c3\u03b5\u03bb\u03af\u03b4\u03b1 \u03b3\u03b9\u03b1 \u03b5\u03b9\u03c3\u03b1\u03b3\u03c9\u03b3\u03ae \u03ba\u03b1\u03b9 \u03b4\u03b9\u03b1\u03c7\u03b5\u03af\u03c1\u03b9\u03c3\u03b7 \" +\n  \"\u03c7\u03c1\u03b7\u03c3\u03c4\u03ce\u03bd/\u03ba\u03b1\u03ba\u03bf\u03c0\u03bf\u03b9\u03ce\u03bd, \u03b4\u03b5\u03bd \u03b8\u03b5\u03c9\u03c1\u03ce \u03cc\u03c4\u03b9 \u03c7\u03c1\u03b5\u03b9\u03b1\u03b6\u03cc\u03bc\u03b1\u03c3\u03c4\u03b5 \u03bd\u03ad\u03b1 \u03c3\u03b5\u03bb\u03af\u03b4\u03b1 \u03b3\u03b9\u03b1 \u03c4\u03b7 \u03b4\u03b9\u03b1\u03b3\u03c1\u03b1\u03c6\u03ae \" +\n  \"\u03c4\u03bf\u03c5\u03c2. \u0398\u03b1 \u03bc\u03c0\u03bf\u03c1\u03bf\u03cd\u03c3\u03b5 \u03bd\u03b1 \u03b5\u03af\u03bd\u03b1\u03b9 \u03b1\u03c0\u03bb\u03ac \u03ad\u03bd\u03b1 \u03ba\u03bf\u03c5\u03bc\u03c0\u03af \u0394\u03b9\u03b1\u03b3\u03c1\u03b1\u03c6\u03ae\u03c2 \u03ba\u03b1\u03b9 \u03bd\u03b1 \u03b4\u03b9\u03b1\u03b3\u03c1\u03ac\u03c6\u03bf\u03c5\u03bc\u03b5 \" +\n  \"\u03b1\u03c0\u03cc \u03c4\u03bf\u03bd \u03c3\u03c5\u03bd\u03bf\u03bb\u03b9\u03ba\u03cc \u03c0\u03af\u03bd\u03b1\u03ba\u03b1.\";\n\nconst firstPara = row2Paragraphs.items[0];\nconst mergedParaXml =\n  \"<w:p><w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t>\" + mergedText + \"</w:t></w:r></w:p>\";\nfirstPara.insertOoxml(wrapOoxml(mergedParaXml), \"Replace\");\nawait context.sync();\n\n// --- Part 2: fill in row 3's (index 2) empty content cell with the new\n// three paragraphs of feedback + fix, ending with the relocated bookmark.\nconst row3Cell = table.getCell(2, 1);\nconst row3Paragraphs = row3Cell.body.paragraphs;\nrow3Paragraphs.load(\"text\");\nawait context.sync();\n\nconst targetPara = row3Paragraphs.items[0];\n\nconst newCellXml =\n  \"<w:p><w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\">\u0394\u03b5\u03bd \u03c5\u03c0\u03ac\u03c1\u03c7\u03b5\u03b9 </w:t></w:r>\" +\n  \"<w:r><w:t>confirmation</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\"> </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\">\u03c3\u03c4\u03bf </w:t></w:r>\" +\n  \"<w:r><w:t>delete</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\">, </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\">\u03bc\u03c0\u03bf\u03c1\u03b5\u03af \u03bd\u03b1 \u03c0\u03b1\u03c4\u03b7\u03b8\u03b5\u03af </w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellStart\\\"/>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t>\u03ba\u03b1\u03c4\u03b1\u03bb\u03ac\u03b8\u03bf\u03c2</w:t></w:r>\" +\n  \"<w:proofErr w:type=\\\"spellEnd\\\"/>\" +\n  \"</w:p>\" +\n  \"<w:p><w:pPr><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr></w:pPr></w:p>\" +\n  \"<w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/>\" +\n  \"<w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"1\\\"/></w:numPr>\" +\n  \"<w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr></w:pPr>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\">\u0394\u03b9\u03bf\u03c1\u03b8\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c0\u03c1\u03bf\u03c3\u03c4\u03ad\u03b8\u03b7\u03ba\u03b5 </w:t></w:r>\" +\n  \"<w:r><w:t>alert</w:t></w:r>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t xml:space=\\\"preserve\\\"> </w:t></w:r>\" +\n  \"<w:r><w:rPr><w:lang w:val=\\\"el-GR\\\"/></w:rPr>\" +\n  \"<w:t>\u03bc\u03ae\u03bd\u03c5\u03bc\u03b1 \u03b3\u03b9\u03b1 \u03c4\u03b7\u03bd \u03b4\u03b9\u03b1\u03b3\u03c1\u03b1\u03c6\u03ae</w:t></w:r>\" +\n  \"<w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/>\" +\n  \"<w:bookmarkEnd w:id=\\\"0\\\"/>\" +\n  \"</w:p>\";\n\ntargetPara.insertOoxml(wrapOoxml(newCellXml), \"Replace\");\nawait context.sync();\n", "ps1": "# Fixed alert confirmation popup not working in DELETE\n#\n# 1) Row 2's first paragraph currently has its Greek text split across two\n#    runs (with a leftover bookmarkStart/bookmarkEnd named \"_GoBack\" sitting\n#    between them, left over from where the cursor last was). Merge that\n#    back into a single contiguous run and drop the stray bookmark.\n# 2) Row 3's content cell is empty. Fill it in with the new feedback (that\n#    there was no confirmation on delete, so it could be clicked by\n#    accident) and the corresponding fix note (an alert/confirm message was\n#    added on delete) - moving the \"_GoBack\" bookmark to sit at the very end\n#    of that new content, which is where the author's cursor ended up.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$mergedParaXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t>\u0395\u03bd\u03ce \u03ad\u03c7\u03b5\u03b9 \u03bd\u03cc\u03b7\u03bc\u03b1 \u03bd\u03b1 \u03ad\u03c7\u03bf\u03c5\u03bc\u03b5 \u03be\u03b5\u03c7\u03c9\u03c1\u03b9\u03c3\u03c4\u03ae \u03c3\u03b5\u03bb\u03af\u03b4\u03b1 \u03b3\u03b9\u03b1 \u03b5\u03b9\u03c3\u03b1\u03b3\u03c9\u03b3\u03ae \u03ba\u03b1\u03b9 \u03b4\u03b9\u03b1\u03c7\u03b5\u03af\u03c1\u03b9\u03c3\u03b7 \u03c7\u03c1\u03b7\u03c3\u03c4\u03ce\u03bd/\u03ba\u03b1\u03ba\u03bf\u03c0\u03bf\u03b9\u03ce\u03bd, \u03b4\u03b5\u03bd \u03b8\u03b5\u03c9\u03c1\u03ce \u03cc\u03c4\u03b9 \u03c7\u03c1\u03b5\u03b9\u03b1\u03b6\u03cc\u03bc\u03b1\u03c3\u03c4\u03b5 \u03bd\u03ad\u03b1 \u03c3\u03b5\u03bb\u03af\u03b4\u03b1 \u03b3\u03b9\u03b1 \u03c4\u03b7 \u03b4\u03b9\u03b1\u03b3\u03c1\u03b1\u03c6\u03ae \u03c4\u03bf\u03c5\u03c2. \u0398\u03b1 \u03bc\u03c0\u03bf\u03c1\u03bf\u03cd\u03c3\u03b5 \u03bd\u03b1 \u03b5\u03af\u03bd\u03b1\u03b9 \u03b1\u03c0\u03bb\u03ac \u03ad\u03bd\u03b1 \u03ba\u03bf\u03c5\u03bc\u03c0\u03af \u0394\u03b9\u03b1\u03b3\u03c1\u03b1\u03c6\u03ae\u03c2 \u03ba\u03b1\u03b9 \u03bd\u03b1 \u03b4\u03b9\u03b1\u03b3\u03c1\u03ac\u03c6\u03bf\u03c5\u03bc\u03b5 \u03b1\u03c0\u03cc \u03c4\u03bf\u03bd \u03c3\u03c5\u03bd\u03bf\u03bb\u03b9\u03ba\u03cc \u03c0\u03af\u03bd\u03b1\u03ba\u03b1.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# --- Part 1: row 2 (\"2\") content cell -> merge the split run, drop bookmark\n$row2Cell = $t.Cell(2, 2)\n$row2FirstPara = $row2Cell.Range.Paragraphs.Item(1)\n$row2FirstPara.Range.InsertXML($mergedParaXml)\n\n$newCellXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t xml:space=\"preserve\">\u0394\u03b5\u03bd \u03c5\u03c0\u03ac\u03c1\u03c7\u03b5\u03b9 </w:t></w:r><w:r><w:t>confirmation</w:t></w:r><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t xml:space=\"preserve\">\u03c3\u03c4\u03bf </w:t></w:r><w:r><w:t>delete</w:t></w:r><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t xml:space=\"preserve\">, </w:t></w:r><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t xml:space=\"preserve\">\u03bc\u03c0\u03bf\u03c1\u03b5\u03af \u03bd\u03b1 \u03c0\u03b1\u03c4\u03b7\u03b8\u03b5\u03af </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t>\u03ba\u03b1\u03c4\u03b1\u03bb\u03ac\u03b8\u03bf\u03c2</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t xml:space=\"preserve\">\u0394\u03b9\u03bf\u03c1\u03b8\u03ce\u03b8\u03b7\u03ba\u03b5, \u03c0\u03c1\u03bf\u03c3\u03c4\u03ad\u03b8\u03b7\u03ba\u03b5 </w:t></w:r><w:r><w:t>alert</w:t></w:r><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:rPr><w:lang w:val=\"el-GR\"/></w:rPr><w:t>\u03bc\u03ae\u03bd\u03c5\u03bc\u03b1 \u03b3\u03b9\u03b1 \u03c4\u03b7\u03bd \u03b4\u03b9\u03b1\u03b3\u03c1\u03b1\u03c6\u03ae</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# --- Part 2: row 3 (\"3\") content cell -> was a single empty paragraph,\n# becomes the three paragraphs above, ending with the relocated bookmark.\n$row3Cell = $t.Cell(3, 2)\n$row3FirstPara = $row3Cell.Range.Paragraphs.Item(1)\n$row3FirstPara.Range.InsertXML($newCellXml)\n\nWrite-Output \"Updated row 2 and row 3 feedback cells.\"\n"}
